# Update "想去人数" (F column) values on the "展览" sheet and the
# combined "全部类型" sheet to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10
$ws1.Range("F7").Value = 1666
$ws1.Range("F8").Value = 17
$ws1.Range("F11").Value = 1599
$ws1.Range("F13").Value = 60
$ws1.Range("F14").Value = 396
$ws1.Range("F21").Value = 180
$ws1.Range("F22").Value = 288
$ws1.Range("F23").Value = 157
$ws1.Range("F25").Value = 222

# --- Sheet "全部类型" (all types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10
$ws4.Range("F7").Value = 1666
$ws4.Range("F9").Value = 17
$ws4.Range("F12").Value = 1599
$ws4.Range("F14").Value = 60
$ws4.Range("F15").Value = 396
$ws4.Range("F22").Value = 180
$ws4.Range("F23").Value = 288
$ws4.Range("F24").Value = 157
$ws4.Range("F26").Value = 222
